# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
#
# Rows 93, 94, 95 (sheet "Venezuela Primera Division") have their
# match-id / result / odds data (columns B, E:AD) cyclically rotated:
#   new row93 <- old row94
#   new row94 <- old row95
#   new row95 <- old row93
# Columns A, C, D stay put for each row (rank index, Div, Date).
#
# Rows 102, 103 simply swap the same set of columns (B, E:AD).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($ws, $row) {
    $data = @{}
    $data.B = $ws.Cells.Item($row, 2).Value2
    for ($c = 5; $c -le 30; $c++) {
        $data[[string]$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $data
}

function Set-RowData($ws, $row, $data) {
    $ws.Cells.Item($row, 2).Value2 = $data.B
    for ($c = 5; $c -le 30; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $data[[string]$c]
    }
}

# --- rotate rows 93 / 94 / 95 ---
$r93 = Get-RowData $ws 93
$r94 = Get-RowData $ws 94
$r95 = Get-RowData $ws 95

Set-RowData $ws 93 $r94
Set-RowData $ws 94 $r95
Set-RowData $ws 95 $r93

# --- swap rows 102 / 103 ---
$r102 = Get-RowData $ws 102
$r103 = Get-RowData $ws 103

Set-RowData $ws 102 $r103
Set-RowData $ws 103 $r102
